{"js": "// Update the date line and all \"A\u00d7B=C\" answer cells to the new values\n// described by the commit. Each old text string is unique in the\n// document, so a simple search-and-replace per pair is safe and keeps\n// the original run formatting (font, size, etc.) untouched.\nconst replacements = [\n  [\"2025-01-08 Wednesday\", \"2025-01-14 Tuesday\"],\n  [\"802\u00d76=4812\", \"867\u00d78=6936\"],\n  [\"103\u00d72=206\", \"222\u00d74=888\"],\n  [\"330\u00d79=2970\", \"364\u00d75=1820\"],\n  [\"679\u00d78=5432\", \"598\u00d74=2392\"],\n  [\"566\u00d75=2830\", \"850\u00d79=7650\"],\n  [\"736\u00d76=4416\", \"651\u00d75=3255\"],\n  [\"622\u00d75=3110\", \"603\u00d77=4221\"],\n  [\"560\u00d74=2240\", \"203\u00d76=1218\"],\n  [\"730\u00d75=3650\", \"623\u00d75=3115\"],\n  [\"716\u00d79=6444\", \"336\u00d76=2016\"],\n  [\"140\u00d78=1120\", \"977\u00d77=6839\"],\n  [\"827\u00d72=1654\", \"593\u00d75=2965\"],\n  [\"641\u00d74=2564\", \"166\u00d75=830\"],\n  [\"463\u00d77=3241\", \"936\u00d74=3744\"],\n  [\"981\u00d76=5886\", \"768\u00d75=3840\"],\n  [\"900\u00d76=5400\", \"869\u00d75=4345\"],\n  [\"407\u00d79=3663\", \"118\u00d79=1062\"],\n  [\"955\u00d76=5730\", \"539\u00d78=4312\"],\n  [\"581\u00d76=3486\", \"988\u00d74=3952\"],\n  [\"359\u00d74=1436\", \"781\u00d78=6248\"],\n  [\"824\u00d72=1648\", \"279\u00d79=2511\"],\n  [\"931\u00d75=4655\", \"766\u00d75=3830\"],\n  [\"298\u00d75=1490\", \"250\u00d73=750\"],\n  [\"214\u00d74=856\", \"610\u00d79=5490\"],\n  [\"892\u00d79=8028\", \"501\u00d78=4008\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all \"A\u00d7B=C\" answer cells to the new values\n# described by the commit. Each old text string is unique in the\n# document, so Find/Replace per pair is safe and keeps the original\n# run formatting (font, size, etc.) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-08 Wednesday\", \"2025-01-14 Tuesday\"),\n    @(\"802\u00d76=4812\", \"867\u00d78=6936\"),\n    @(\"103\u00d72=206\", \"222\u00d74=888\"),\n    @(\"330\u00d79=2970\", \"364\u00d75=1820\"),\n    @(\"679\u00d78=5432\", \"598\u00d74=2392\"),\n    @(\"566\u00d75=2830\", \"850\u00d79=7650\"),\n    @(\"736\u00d76=4416\", \"651\u00d75=3255\"),\n    @(\"622\u00d75=3110\", \"603\u00d77=4221\"),\n    @(\"560\u00d74=2240\", \"203\u00d76=1218\"),\n    @(\"730\u00d75=3650\", \"623\u00d75=3115\"),\n    @(\"716\u00d79=6444\", \"336\u00d76=2016\"),\n    @(\"140\u00d78=1120\", \"977\u00d77=6839\"),\n    @(\"827\u00d72=1654\", \"593\u00d75=2965\"),\n    @(\"641\u00d74=2564\", \"166\u00d75=830\"),\n    @(\"463\u00d77=3241\", \"936\u00d74=3744\"),\n    @(\"981\u00d76=5886\", \"768\u00d75=3840\"),\n    @(\"900\u00d76=5400\", \"869\u00d75=4345\"),\n    @(\"407\u00d79=3663\", \"118\u00d79=1062\"),\n    @(\"955\u00d76=5730\", \"539\u00d78=4312\"),\n    @(\"581\u00d76=3486\", \"988\u00d74=3952\"),\n    @(\"359\u00d74=1436\", \"781\u00d78=6248\"),\n    @(\"824\u00d72=1648\", \"279\u00d79=2511\"),\n    @(\"931\u00d75=4655\", \"766\u00d75=3830\"),\n    @(\"298\u00d75=1490\", \"250\u00d73=750\"),\n    @(\"214\u00d74=856\", \"610\u00d79=5490\"),\n    @(\"892\u00d79=8028\", \"501\u00d78=4008\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
